$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells that will hold numeric-looking text to stay as text
$ws.Range("D2:D51").NumberFormat = "@"

# Update Price (D) column values
$ws.Range("D2").Value = "69.075.75"
$ws.Range("D3").Value = "3.520.13"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "578.43"
$ws.Range("D6").Value = "171.34"
$ws.Range("D7").Value = "0.620"
$ws.Range("D8").Value = "3.511.06"
$ws.Range("D9").Value = "1.00"
$ws.Range("D11").Value = "6.70"
$ws.Range("D12").Value = "0.600"
$ws.Range("D13").Value = "47.29"
$ws.Range("D15").Value = "685.36"
$ws.Range("D16").Value = "4.083.30"
$ws.Range("D17").Value = "8.80"
$ws.Range("D18").Value = "69.140.89"
$ws.Range("D19").Value = "3.513.28"
$ws.Range("D21").Value = "17.42"
$ws.Range("D22").Value = "11.17"
$ws.Range("D25").Value = "97.65"
$ws.Range("D30").Value = "33.28"
$ws.Range("D33").Value = "1.36"
$ws.Range("D35").Value = "571.92"
$ws.Range("D36").Value = "3.65"
$ws.Range("D39").Value = "57.12"
$ws.Range("D40").Value = "1.00"
$ws.Range("D42").Value = "0.0439"
$ws.Range("D44").Value = "3.435.59"
$ws.Range("D45").Value = "33.25"
$ws.Range("D46").Value = "0.0₃0703"
$ws.Range("D49").Value = "0.134"
$ws.Range("D50").Value = "134.17"
$ws.Range("D51").Value = "0.150"

# Update Volume(1h) (E) column values
$ws.Range("E2").Value = "  -3.88%  "
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("E6").Value = "  -5.62%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -6.45%  "
$ws.Range("E11").Value = "  +13.32%  "
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("E13").Value = "  -4.87%  "
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("E18").Value = "  -3.82%  "
$ws.Range("E19").Value = "  -3.89%  "
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  -4.81%  "
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("E24").Value = "  -6.93%  "
$ws.Range("E25").Value = "  -5.37%  "
$ws.Range("E26").Value = "  -4.29%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -6.52%  "
$ws.Range("E29").Value = "  -5.77%  "
$ws.Range("E30").Value = "  -5.53%  "
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("E32").Value = "  -7.04%  "
$ws.Range("E33").Value = "  -5.61%  "
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").Value = "  -2.22%  "
$ws.Range("E36").Value = "  -13.90%  "
$ws.Range("E37").Value = "  -4.14%  "
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("E42").Value = "  -6.14%  "
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("E44").Value = "  -6.49%  "
$ws.Range("E45").Value = "  -6.91%  "
$ws.Range("E46").Value = "  -7.82%  "
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").Value = "  -7.10%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("E51").Value = "  -0.21%  "
